$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing text storage
# (matches the source data which keeps numeric-looking price/volume
# strings as literal text, e.g. "26.311.30", "238.83", "  +2.88%  ").
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '26.311.30'
Set-TextCell 'E2' '  +2.88%  '

Set-TextCell 'D3' '1.718.11'
Set-TextCell 'E3' '  +3.07%  '

Set-TextCell 'E4' '  +0.05%  '

Set-TextCell 'D5' '238.83'
Set-TextCell 'E5' '  +0.83%  '

Set-TextCell 'E6' '  +0.06%  '

Set-TextCell 'D7' '0.4714'
Set-TextCell 'E7' '  -1.92%  '

Set-TextCell 'D8' '0.2623'
Set-TextCell 'E8' '  -0.54%  '

Set-TextCell 'D9' '0.06182'
Set-TextCell 'E9' '  +0.44%  '

Set-TextCell 'D10' '1.718.09'
Set-TextCell 'E10' '  +2.97%  '

Set-TextCell 'E11' '  -0.43%  '

Set-TextCell 'D12' '15.32'
Set-TextCell 'E12' '  +3.14%  '

Set-TextCell 'D13' '0.5912'
Set-TextCell 'E13' '  -1.31%  '

Set-TextCell 'D14' '4.399'
Set-TextCell 'E14' '  -0.25%  '

Set-TextCell 'D15' '76.11'
Set-TextCell 'E15' '  +2.10%  '

Set-TextCell 'E16' '  +0.01%  '

Set-TextCell 'E17' '  +0.11%  '

Set-TextCell 'D18' '26.307.55'

Set-TextCell 'D19' '0.000006788'
Set-TextCell 'E19' '  -0.04%  '

Set-TextCell 'E20' '  +0.70%  '

Set-TextCell 'D21' '1.938.23'
Set-TextCell 'E21' '  +3.15%  '

Set-TextCell 'D22' '4.542'
Set-TextCell 'E22' '  +1.61%  '

Set-TextCell 'D23' '8.728'
Set-TextCell 'E23' '  +0.21%  '

Set-TextCell 'D24' '5.321'
Set-TextCell 'E24' '  -0.65%  '

Set-TextCell 'D25' '135.83'
Set-TextCell 'E25' '  +1.05%  '

Set-TextCell 'D26' '15.23'
Set-TextCell 'E26' '  +0.79%  '

Set-TextCell 'D27' '1.409'
Set-TextCell 'E27' '  +0.21%  '

Set-TextCell 'D28' '108.21'
Set-TextCell 'E28' '  +3.05%  '

Set-TextCell 'D29' '1.755'
Set-TextCell 'E29' '  +3.52%  '

Set-TextCell 'D30' '3.998'

Set-TextCell 'E31' '  +0.40%  '

Set-TextCell 'E32' '  +0.47%  '

Set-TextCell 'D33' '0.04433'
Set-TextCell 'E33' '  +1.71%  '

Set-TextCell 'D34' '2.614'
Set-TextCell 'E34' '  -0.09%  '

Set-TextCell 'D35' '0.9748'
Set-TextCell 'E35' '  +2.37%  '

Set-TextCell 'D36' '0.6187'
Set-TextCell 'E36' '  +0.02%  '

Set-TextCell 'D37' '0.9269'
Set-TextCell 'E37' '  +6.39%  '

Set-TextCell 'D38' '114.11'
Set-TextCell 'E38' '  +16.59%  '

Set-TextCell 'D39' '2.414'
Set-TextCell 'E39' '  -7.64%  '

Set-TextCell 'D40' '1.001'
Set-TextCell 'E40' '  +0.02%  '

Set-TextCell 'D41' '1.896'
Set-TextCell 'E41' '  +0.66%  '

Set-TextCell 'D42' '0.01480'
Set-TextCell 'E42' '  -2.36%  '

Set-TextCell 'D43' '5.333'
Set-TextCell 'E43' '  +13.71%  '

Set-TextCell 'D44' '0.3806'
Set-TextCell 'E44' '  +0.68%  '

Set-TextCell 'D45' '0.1159'
Set-TextCell 'E45' '  +3.07%  '

Set-TextCell 'D46' '6.274'
Set-TextCell 'E46' '  +0.66%  '

Set-TextCell 'D47' '0.05289'
Set-TextCell 'E47' '  +0.50%  '

Set-TextCell 'D48' '30.47'
Set-TextCell 'E48' '  +3.16%  '

Set-TextCell 'D49' '7.668'
Set-TextCell 'E49' '  +3.51%  '

# Rows 50/51: NEARProtocol and Decentraland swap rank positions, each
# with refreshed price/volume figures.
Set-TextCell 'B50' 'NEARProtocol'
Set-TextCell 'C50' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D50' '1.215'
Set-TextCell 'E50' '  +1.52%  '

Set-TextCell 'B51' 'Decentraland'
Set-TextCell 'C51' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 'D51' '0.3368'
Set-TextCell 'E51' '  +0.62%  '
